$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2:F4").Value = "NO"
$ws.Range("G2:G4").Value = "user already created"
